$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-22 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-23 Monday", 2)

$d.Content.Find.Execute("694÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "406÷4=", 2)
$d.Content.Find.Execute("280÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "410÷6=", 2)
$d.Content.Find.Execute("225÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "978÷8=", 2)
$d.Content.Find.Execute("519÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "848÷3=", 2)
$d.Content.Find.Execute("710÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "815÷8=", 2)
$d.Content.Find.Execute("732÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "331÷7=", 2)
$d.Content.Find.Execute("852÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "997÷3=", 2)
$d.Content.Find.Execute("733÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "195÷2=", 2)
$d.Content.Find.Execute("575÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "962÷5=", 2)
$d.Content.Find.Execute("425÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "508÷7=", 2)
$d.Content.Find.Execute("726÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "480÷7=", 2)
$d.Content.Find.Execute("687÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "831÷6=", 2)
$d.Content.Find.Execute("952÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "522÷3=", 2)
$d.Content.Find.Execute("647÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "272÷2=", 2)
$d.Content.Find.Execute("194÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "906÷2=", 2)
$d.Content.Find.Execute("315÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "648÷3=", 2)
$d.Content.Find.Execute("620÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "470÷5=", 2)
$d.Content.Find.Execute("267÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "872÷8=", 2)
$d.Content.Find.Execute("757÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "989÷6=", 2)
$d.Content.Find.Execute("322÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "828÷5=", 2)
$d.Content.Find.Execute("884÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "543÷5=", 2)
$d.Content.Find.Execute("349÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "178÷6=", 2)
$d.Content.Find.Execute("572÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "438÷6=", 2)
$d.Content.Find.Execute("742÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "725÷3=", 2)
$d.Content.Find.Execute("488÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "587÷8=", 2)
